# Actualización automática 2025-09-22 08:22:24
#
# Updates the per-advisor sales figures for "PIEDRA SINTERIZADA",
# "PORCELANATO" and "SAL SOLUBLE" on the "VENTAS POR GRUPO" sheet, the
# "septiembre" column on the "VENTA MENSUAL" sheet, and the derived
# totals/ratios on the "CUMPLIMIENTO MENSUAL" sheet, together with the
# small column-width tweak that Excel recorded alongside the data refresh.

$wb = $excel.ActiveWorkbook

# ── Sheet 1: "VENTAS POR GRUPO" ─────────────────────────────────────────
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# FUENTES PAREDES MARIA FERNANDA / PIEDRA SINTERIZADA
$wsGrupo.Range("L12").Value = 537.34

# MANCHENO PINO HERVIN SANTIAGO / PORCELANATO y SAL SOLUBLE
$wsGrupo.Range("M20").Value = 1874.48
$wsGrupo.Range("O20").Value = 323.48

# Contador de filas con datos para la columna PIEDRA SINTERIZADA
$wsGrupo.Range("L35").Value = "2 de 33"

# ── Sheet 2: "VENTA MENSUAL" ────────────────────────────────────────────
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# FUENTES PAREDES MARIA FERNANDA / septiembre
$wsMensual.Range("F12").Value = 998.71

# MANCHENO PINO HERVIN SANTIAGO / septiembre
$wsMensual.Range("F20").Value = 2633.74

# Total septiembre
$wsMensual.Range("F35").Value = 16567.13

# Columna "julio" (4ta columna) ligeramente más ancha tras la actualización
$wsMensual.Columns.Item(4).ColumnWidth = 13.14

# ── Sheet 3: "CUMPLIMIENTO MENSUAL" ─────────────────────────────────────
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# PIEDRA SINTERIZADA
$wsCumpl.Range("D11").Value = 3761.38
$wsCumpl.Range("E11").Value = -839.1554181472602
$wsCumpl.Range("F11").Value = 1.287163219198992

# PORCELANATO
$wsCumpl.Range("D12").Value = 10084.97
$wsCumpl.Range("E12").Value = 12348.7853751766
$wsCumpl.Range("F12").Value = 0.449544440123441

# SAL SOLUBLE
$wsCumpl.Range("D14").Value = 323.48
$wsCumpl.Range("E14").Value = 1101.4862010375
$wsCumpl.Range("F14").Value = 0.2270088930982912

# TOTAL
$wsCumpl.Range("D15").Value = 16824.36
$wsCumpl.Range("E15").Value = 21918.65881339592
$wsCumpl.Range("F15").Value = 0.4342552675369412
